# Error Calculations and Plots
# Two rows of the "missing data" sample set are dropped ("RM 232" at row 26
# and "SC 92" at row 28), which shifts every subsequent row up by one/two
# positions respectively. A few cells in column D (the randomly-missing
# column) also move to different rows as part of the new missing-data
# pattern for this seed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (original row 26).
$ws.Rows(26).Delete()

# Remove the "SC 92" row - after the delete above it now sits at row 27.
$ws.Rows(27).Delete()

# Column D "missing" pattern shifts for three of the remaining rows.
$ws.Range("D26").Value = ""
$ws.Range("D27").Value = -14.6
$ws.Range("D29").Value = ""
